$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the business-name field: abbreviation -> full Vietnamese term
$ws.Range("C2").Value = "Hộ Kinh Doanh"

# The QR tool filled in previously-blank fields on row 2 (note, phone,
# branch, cccd, customerCode). Phone and cccd are long digit strings with
# a leading zero / no natural numeric meaning, so pre-format those two
# cells as Text before typing them in (otherwise Excel's normal "looks
# like a number" auto-conversion would strip the leading zero / turn them
# into numeric cells). The format is reset back to Normal immediately
# after so the cell keeps the workbook's default (General) style, just
# like the other string cells on the row.
$ws.Range("H2").NumberFormat = "@"
$ws.Range("J2").NumberFormat = "@"

$ws.Range("G2").Value = "Nhà mặt tiền"
$ws.Range("H2").Value = "01234567891"
$ws.Range("I2").Value = "Cần Thơ II"
$ws.Range("J2").Value = "123456789101"
$ws.Range("K2").Value = "1890-000000000"

$ws.Range("H2").Style = "Normal"
$ws.Range("J2").Style = "Normal"
